$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Yellow highlight colour used for all the newly-filled cells (matches the
# workbook's new fill: solid FFFFFF00). COM Color values are 0xBBGGRR.
$yellow = 65535

# --- Fill in the missing symmetric entries in the big comparison table ---
$ws.Range("H5").Value = "SBS"
$ws.Range("H5").Interior.Color = $yellow

$ws.Range("H6").Value = "SBS"
$ws.Range("H6").Interior.Color = $yellow

$ws.Range("H8").Value = "FBS"
$ws.Range("H8").Interior.Color = $yellow

$ws.Range("H9").Value = "FBS"
$ws.Range("H9").Interior.Color = $yellow

$ws.Range("I10").Value = "SAF"
$ws.Range("I10").Interior.Color = $yellow

$ws.Range("I11").Value = "SAF"
$ws.Range("I11").Interior.Color = $yellow

$ws.Range("I13").Value = "FAF"
$ws.Range("I13").Interior.Color = $yellow

$ws.Range("I14").Value = "FAF"
$ws.Range("I14").Interior.Color = $yellow

# --- Add the new "Auto-consequences" notes block below the table ---
$ws.Range("B20").Value = "Auto-consequences:"
$ws.Range("B20").Interior.Color = $yellow

$ws.Range("B21").Value = "SBS > SBF"
$ws.Range("B21").Interior.Color = $yellow

$ws.Range("B22").Value = "FBS > SBS, SBF, FBF"
$ws.Range("B22").Interior.Color = $yellow

$ws.Range("B23").Value = "SAF > FAF, FAS, SAS"
$ws.Range("B23").Interior.Color = $yellow

$ws.Range("B24").Value = "FAF > FAS"
$ws.Range("B24").Interior.Color = $yellow

# --- Update the view state to match: scrolled down one row, new block selected ---
$ws.Activate()
$ws.Range("B21:B24").Select()
